# Swap the deck's theme color palette: theme1.xml (the slide master's
# theme, currently the "Integral" theme) is recolored to the stock
# "Office Theme" palette, mirroring the content that theme2.xml (the
# notes master's theme) already has.
#
# PowerPoint's ThemeColorScheme exposes exactly 12 slots, in this
# fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1 .. 10 accent6,
#   11 hlink, 12 folHlink
# Each slot's .RGB is a VBA-style packed integer: R + G*256 + B*65536
# (i.e. 0x00BBGGRR), not the usual 0xRRGGBB, so build it from the hex
# triples instead of pasting hex literals directly.

function ToVbaRgb([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the standard Office theme colours.
$officeTheme = @(
    @(0x00, 0x00, 0x00),  # 1  dk1
    @(0xFF, 0xFF, 0xFF),  # 2  lt1
    @(0x44, 0x54, 0x6A),  # 3  dk2
    @(0xE7, 0xE6, 0xE6),  # 4  lt2
    @(0x5B, 0x9B, 0xD5),  # 5  accent1
    @(0xED, 0x7D, 0x31),  # 6  accent2
    @(0xA5, 0xA5, 0xA5),  # 7  accent3
    @(0xFF, 0xC0, 0x00),  # 8  accent4
    @(0x44, 0x72, 0xC4),  # 9  accent5
    @(0x70, 0xAD, 0x47),  # 10 accent6
    @(0x05, 0x63, 0xC1),  # 11 hlink
    @(0x95, 0x4F, 0x72)   # 12 folHlink
)

$p = $ppt.ActivePresentation

# Any slide anchors onto the (single) slide master's theme; use slide 1.
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $rgbVal = ToVbaRgb $officeTheme[$i - 1][0] $officeTheme[$i - 1][1] $officeTheme[$i - 1][2]
    $tcs.Colors($i).RGB = $rgbVal
}
